$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitoramento diário")

$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 20
$ws.Range("G7").Value = 20
$ws.Range("I7").Value = 20

$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 14
$ws.Range("F10").Value = 14
$ws.Range("G10").Value = 14
$ws.Range("H10").Value = 14
$ws.Range("I10").Value = 14
$ws.Range("J10").Value = 14

$ws.Range("C13").Value = 2
$ws.Range("E13").Value = 14
$ws.Range("G13").Value = 14
$ws.Range("I13").Value = 14
$ws.Range("J13").Value = 14

$ws.Range("C14").ClearContents()

$ws.Range("C14").Select()
